$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 283, shifting the existing rows 283:312 down to 286:315.
$ws.Rows("283:285").Insert()

# New row 283: Acelga, Extra, week of 2021-09-22 (serial 44461)
$ws.Cells.Item(283, 1).Value = 9
$ws.Cells.Item(283, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(283, 3).Value = "Metropolitana"
$ws.Cells.Item(283, 4).Value = 44461
$ws.Cells.Item(283, 5).Value = 13
$ws.Cells.Item(283, 6).Value = 100112009
$ws.Cells.Item(283, 7).Value = "Acelga"
$ws.Cells.Item(283, 8).Value = "Sin especificar"
$ws.Cells.Item(283, 9).Value = "Extra"
$ws.Cells.Item(283, 10).Value = 16
$ws.Cells.Item(283, 11).Value = 11000
$ws.Cells.Item(283, 12).Value = 12000
$ws.Cells.Item(283, 13).Value = 11500
$ws.Cells.Item(283, 14).Value = "$/docena de atados"
$ws.Cells.Item(283, 15).Value = "Región Metropolitana"
$ws.Cells.Item(283, 16).Value = 3833
$ws.Cells.Item(283, 17).Value = 3
$ws.Cells.Item(283, 18).Value = "Hortaliza"

# New row 284: Acelga, Primera, week of 2021-09-22 (serial 44461)
$ws.Cells.Item(284, 1).Value = 9
$ws.Cells.Item(284, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(284, 3).Value = "Metropolitana"
$ws.Cells.Item(284, 4).Value = 44461
$ws.Cells.Item(284, 5).Value = 13
$ws.Cells.Item(284, 6).Value = 100112009
$ws.Cells.Item(284, 7).Value = "Acelga"
$ws.Cells.Item(284, 8).Value = "Sin especificar"
$ws.Cells.Item(284, 9).Value = "Primera"
$ws.Cells.Item(284, 10).Value = 34
$ws.Cells.Item(284, 11).Value = 9000
$ws.Cells.Item(284, 12).Value = 10000
$ws.Cells.Item(284, 13).Value = 9500
$ws.Cells.Item(284, 14).Value = "$/docena de atados"
$ws.Cells.Item(284, 15).Value = "Región Metropolitana"
$ws.Cells.Item(284, 16).Value = 3167
$ws.Cells.Item(284, 17).Value = 3
$ws.Cells.Item(284, 18).Value = "Hortaliza"

# New row 285: Acelga, Segunda, week of 2021-09-22 (serial 44461)
$ws.Cells.Item(285, 1).Value = 9
$ws.Cells.Item(285, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(285, 3).Value = "Metropolitana"
$ws.Cells.Item(285, 4).Value = 44461
$ws.Cells.Item(285, 5).Value = 13
$ws.Cells.Item(285, 6).Value = 100112009
$ws.Cells.Item(285, 7).Value = "Acelga"
$ws.Cells.Item(285, 8).Value = "Sin especificar"
$ws.Cells.Item(285, 9).Value = "Segunda"
$ws.Cells.Item(285, 10).Value = 25
$ws.Cells.Item(285, 11).Value = 7000
$ws.Cells.Item(285, 12).Value = 8000
$ws.Cells.Item(285, 13).Value = 7480
$ws.Cells.Item(285, 14).Value = "$/docena de atados"
$ws.Cells.Item(285, 15).Value = "Región Metropolitana"
$ws.Cells.Item(285, 16).Value = 2493
$ws.Cells.Item(285, 17).Value = 3
$ws.Cells.Item(285, 18).Value = "Hortaliza"
